$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

$ws.Range("B13").Value = "01/01/2023"
$ws.Range("C13").Value = "01/01/2023"

$ws.Range("B15").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C15").Value = "519033 - Carlos Yujiro Shigue"

$ws.Range("B18").Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("C18").Value = "7290967 - Emerson Gonçalves de Melo"
